# Generate Report for Handback
#
# Renames the two tracked files in the handback-status report:
#   7678ce99-c951-4420-978b-50240a635986  ->  05049e00-37bb-4c68-ae9a-126150ae4e7f
#   82bf6c72-145b-4e01-8ed8-41f688a28fec  ->  ffff0c29e0a4-cc9d-4261-bd52-5c0058008b0d
# and refreshes the associated handoff/handback timestamps + generated xliff
# file names across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "7678ce99-c951-4420-978b-50240a635986"
$newGuid1 = "05049e00-37bb-4c68-ae9a-126150ae4e7f"
$oldGuid2 = "82bf6c72-145b-4e01-8ed8-41f688a28fec"
$newGuid2 = "ffff0c29e0a4-cc9d-4261-bd52-5c0058008b0d"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newGuid1.md"
$ws1.Range("B2").Value = "e2e\$newGuid1.md"
$ws1.Range("G2").Value = "2016-08-15 14:59:51"

$ws1.Range("A3").Value = "$newGuid2.md"
$ws1.Range("B3").Value = "e2e\$newGuid2.md"
$ws1.Range("G3").Value = "2016-08-15 14:59:51"

foreach ($hl in $ws1.Hyperlinks) {
    if ($hl.TextToDisplay -eq "e2e\$oldGuid1.md") {
        $hl.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($hl.TextToDisplay -eq "e2e\$oldGuid2.md") {
        $hl.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newGuid1.md"
$ws2.Range("G2").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-15 14:59:44"
$ws2.Range("I2").Value = "$newGuid1.md"
$ws2.Range("J2").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-15 15:00:34"

$ws2.Range("A3").Value = "$newGuid2.md"
$ws2.Range("G3").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-15 14:59:44"
$ws2.Range("I3").Value = "$newGuid2.md"
$ws2.Range("J3").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-15 15:00:34"

foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid1.md") {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($hl.TextToDisplay -eq "$oldGuid2.md") {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newGuid1.md"
$ws3.Range("G2").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-15 14:59:51"
$ws3.Range("I2").Value = "$newGuid1.md"
$ws3.Range("J2").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-15 15:00:42"

$ws3.Range("A3").Value = "$newGuid2.md"
$ws3.Range("G3").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-15 14:59:51"
$ws3.Range("I3").Value = "$newGuid2.md"
$ws3.Range("J3").Value = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-15 15:00:42"

foreach ($hl in $ws3.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid1.md") {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($hl.TextToDisplay -eq "$oldGuid2.md") {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}
